$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (test case 5 "Logga ut"): the "Handling" instructions now refer to the
# "Inloggad som" menu instead of "Användarsida".
$ws.Range("C6").Value = "Välj `"Inloggad som`" och tryck `"Logga ut`"."

# Row 8 (test case 7 "Placera en order"): the "Resultat" column was a stray
# "Null" placeholder, now a real expected result.
$ws.Range("E8").Value = "Ordern är betald"

# The "Not" column notes for rows 7 & 8 ("Backend" / "Ej klar") are done, so
# clear them (this also drops their Neutral/Dålig highlight styling below).
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()

# Row 9 (test case 8 "Lämna en recension"): note updated from "Backend" to
# "Back end klar".
$ws.Range("G9").Value = "Back end klar"

# F7 was highlighted "Neutral" (yellow) and F8 was highlighted "Dålig" (red)
# as reminders; both go back to the plain Normal style now that the sprint
# is being handed in.
$ws.Range("F7").Style = "Normal"
$ws.Range("F8").Style = "Normal"

# With no cell referencing the "Dålig" cell style any more, drop it from the
# workbook's style list entirely.
$styles = $wb.Styles()
$badStyle = $styles.Item("Dålig")
$badStyle.Delete()

# Column G was widened a bit.
$ws.Columns.Item(7).ColumnWidth = 11.85

# Update the saved cursor/selection position.
$ws.Range("I13").Select()
